$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2300
$ws.Range("J40").Value = 2300
$ws.Range("L40").Value = 2300
$ws.Range("N40").Value = -2650

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2384.111
$ws.Range("I61").Value = 2083.0908
$ws.Range("K61").Value = 2083.0908
$ws.Range("M61").Value = -1871.0908
$ws.Range("H62").Value = 14249
$ws.Range("J62").Value = 14249
$ws.Range("L62").Value = 14249
$ws.Range("N62").Value = -15497
$ws.Range("H65").Value = 14249
$ws.Range("J65").Value = 14249
$ws.Range("L65").Value = 42747
$ws.Range("N65").Value = -48987
$ws.Range("H76").Value = 33155.2
$ws.Range("I76").Value = 24000
$ws.Range("K76").Value = 24000
$ws.Range("M76").Value = -23662
$ws.Range("H79").Value = 33155.2
$ws.Range("I79").Value = 24000
$ws.Range("K79").Value = 24000
$ws.Range("M79").Value = -22830
$ws.Range("H80").Value = 25957.4
$ws.Range("J80").Value = 25957.4
$ws.Range("L80").Value = 25957.4
$ws.Range("N80").Value = -27953.4
$ws.Range("H83").Value = 25957.4
$ws.Range("J83").Value = 25957.4
$ws.Range("L83").Value = 77872.20000000001
$ws.Range("N83").Value = -87856.20000000001
$ws.Range("H129").Value = 48867.832
$ws.Range("I129").Value = 47709
$ws.Range("K129").Value = 47709
$ws.Range("M129").Value = -42709
$ws.Range("H132").Value = 4140.2666
$ws.Range("I132").Value = 4865.778
$ws.Range("J132").Value = 3052
$ws.Range("K132").Value = 14597.334
$ws.Range("L132").Value = 9156
$ws.Range("M132").Value = -12067.334
$ws.Range("N132").Value = -14216
$ws.Range("H136").Value = 2384.111
$ws.Range("I136").Value = 2083.0908
$ws.Range("K136").Value = 6249.2724
$ws.Range("M136").Value = -3699.2724

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 34757.6
$ws.Range("J35").Value = 34757.6
$ws.Range("L35").Value = 34757.6
$ws.Range("N35").Value = -35377.6
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H134").Value = 1649.8334
$ws.Range("I134").Value = 1458.909
$ws.Range("K134").Value = 4376.727000000001
$ws.Range("M134").Value = -1841.727000000001
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 9728
$ws.Range("J50").Value = 11073.6
$ws.Range("L50").Value = 11073.6
$ws.Range("N50").Value = -12323.6
$ws.Range("H51").Value = 9362.5
$ws.Range("J51").Value = 10983.333
$ws.Range("L51").Value = 10983.333
$ws.Range("N51").Value = -12455.333
$ws.Range("H58").Value = 1973
$ws.Range("I58").Value = 1568.2222
$ws.Range("J58").Value = 2377.7778
$ws.Range("K58").Value = 1568.2222
$ws.Range("L58").Value = 2377.7778
$ws.Range("M58").Value = -1365.2222
$ws.Range("N58").Value = -2783.7778
$ws.Range("H60").Value = 13876.375
$ws.Range("J60").Value = 13876.375
$ws.Range("L60").Value = 13876.375
$ws.Range("N60").Value = -14898.375
$ws.Range("H61").Value = 9362.5
$ws.Range("J61").Value = 10983.333
$ws.Range("L61").Value = 10983.333
$ws.Range("N61").Value = -11679.333
$ws.Range("H109").Value = 34745.2
$ws.Range("J109").Value = 34745.2
$ws.Range("L109").Value = 34745.2
$ws.Range("N109").Value = -36825.2
$ws.Range("H132").Value = 2827.12
$ws.Range("I132").Value = 2486.5881
$ws.Range("J132").Value = 3550.75
$ws.Range("K132").Value = 7459.7643
$ws.Range("L132").Value = 10652.25
$ws.Range("M132").Value = -4929.7643
$ws.Range("N132").Value = -15712.25
$ws.Range("H136").Value = 1973
$ws.Range("I136").Value = 1568.2222
$ws.Range("J136").Value = 2377.7778
$ws.Range("K136").Value = 4704.6666
$ws.Range("L136").Value = 7133.3334
$ws.Range("M136").Value = -2154.6666
$ws.Range("N136").Value = -12233.3334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 16106.692
$ws.Range("I57").Value = 980
$ws.Range("K57").Value = 980
$ws.Range("M57").Value = -160
$ws.Range("H103").Value = 26000
$ws.Range("J103").Value = 26000
$ws.Range("L103").Value = 26000
$ws.Range("N103").Value = -28344
$ws.Range("H123").Value = 32995
$ws.Range("J123").Value = 32995
$ws.Range("L123").Value = 32995
$ws.Range("N123").Value = -37895
$ws.Range("H126").Value = 3613.6858
$ws.Range("I126").Value = 1946.9131
$ws.Range("K126").Value = 5840.7393
$ws.Range("M126").Value = -3370.7393
$ws.Range("H132").Value = 2601.7693
$ws.Range("I132").Value = 2120.111
$ws.Range("J132").Value = 3685.5
$ws.Range("K132").Value = 6360.333
$ws.Range("L132").Value = 11056.5
$ws.Range("M132").Value = -3830.333
$ws.Range("N132").Value = -16116.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 21125
$ws.Range("J109").Value = 21125
$ws.Range("L109").Value = 21125
$ws.Range("N109").Value = -23899
$ws.Range("H132").Value = 2946.6924
$ws.Range("I132").Value = 2125.5
$ws.Range("J132").Value = 4260.6
$ws.Range("K132").Value = 6376.5
$ws.Range("L132").Value = 12781.8
$ws.Range("M132").Value = -3846.5
$ws.Range("N132").Value = -17841.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("H29").Value = 18330
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 18330
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 18330
$ws.Range("N29").Value = -18910
$ws.Range("H109").Value = 24394.25
$ws.Range("J109").Value = 24394.25
$ws.Range("L109").Value = 24394.25
$ws.Range("N109").Value = -27168.25
$ws.Range("H132").Value = 2418.7896
$ws.Range("I132").Value = 2293.88
$ws.Range("J132").Value = 2659
$ws.Range("K132").Value = 6881.64
$ws.Range("L132").Value = 7977
$ws.Range("M132").Value = -4351.64
$ws.Range("N132").Value = -13037
$ws.Range("H133").Value = 39476.668
$ws.Range("J133").Value = 39476.668
$ws.Range("L133").Value = 39476.668
$ws.Range("N133").Value = -49596.668
$ws.Range("H136").Value = 1211.1731
$ws.Range("I136").Value = 1129.2565
$ws.Range("J136").Value = 1456.9231
$ws.Range("K136").Value = 3387.7695
$ws.Range("L136").Value = 4370.7693
$ws.Range("M136").Value = -837.7694999999999
$ws.Range("N136").Value = -9470.7693
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("M29").ClearContents()
